# Refresh cryptos price (D) and volume-change (E) columns to match latest scrape.
# Values that look like plain numbers are entered with a leading apostrophe so
# Excel keeps them as text (matching the original inline/shared-string cells)
# instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.971.87'
$ws.Range("D3").Value = '1.562.06'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'207.35"
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = "'22.11"
$ws.Range("E8").Value = '  +0.91%  '
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("E10").Value = '  +2.65%  '
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").Value = '1.785.79'
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D13").Value = '1.563.29'
$ws.Range("E13").Value = '  +0.62%  '
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("D15").Value = "'0.519"
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").Value = "'61.94"
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").Value = '26.959.48'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("D19").Value = "'215.74"
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").Value = "'7.36"
$ws.Range("E20").Value = '  +1.10%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = "'4.10"
$ws.Range("E22").Value = '  +1.56%  '
$ws.Range("D23").Value = "'9.18"
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = "'1.92"
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("D25").Value = "'153.24"
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("D26").Value = "'6.61"
$ws.Range("E26").Value = '  +0.43%  '
$ws.Range("E27").Value = '  +1.23%  '
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = "'0.0470"
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("E31").Value = '  +1.79%  '
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("D33").Value = "'3.11"
$ws.Range("E33").Value = '  +1.66%  '
$ws.Range("D34").Value = '1.420.51'
$ws.Range("E34").Value = '  -1.07%  '
$ws.Range("D35").Value = "'1.60"
$ws.Range("E35").Value = '  +2.76%  '
$ws.Range("D36").Value = "'1.07"
$ws.Range("E36").Value = '  +9.36%  '
$ws.Range("E37").Value = '  +2.37%  '
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = "'0.535"
$ws.Range("E39").Value = '  +3.08%  '
$ws.Range("E40").Value = '  +2.39%  '
$ws.Range("D41").Value = "'0.807"
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("E43").Value = '  +2.46%  '
$ws.Range("E44").Value = '  +2.05%  '
$ws.Range("D45").Value = "'64.56"
$ws.Range("E45").Value = '  +0.93%  '
$ws.Range("D46").Value = "'1.74"
$ws.Range("E46").Value = '  -1.03%  '
$ws.Range("D47").Value = '1.698.82'
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("D48").Value = "'87.20"
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("D50").Value = '0.0₇0999'
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("D51").Value = "'0.0958"
$ws.Range("E51").Value = '  +0.23%  '
